$d = $word.ActiveDocument

# Locate the target paragraph by searching for the start of its text,
# then expand the found range to cover the whole paragraph (minus the
# trailing paragraph mark).
$findRng = $d.Content
$findRng.Find.ClearFormatting()
$ok = $findRng.Find.Execute( `
    "The data suggests that there is more variability with successful campaigns", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $ok) {
    throw "Could not locate target paragraph text"
}

$null = $findRng.Expand(4)  # wdParagraph - expand to the whole paragraph, incl. the mark

$paraStart = $findRng.Start
$paraEnd = $findRng.End - 1   # exclude the trailing paragraph mark

$targetRng = $d.Range($paraStart, $paraEnd)

# Build an OOXML fragment (flat WordprocessingML package) that splits the
# paragraph's single run into three runs, each carrying an explicit Arial
# rFonts rPr, and with the comma separated from its neighbouring text.
$xmlFragment = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
'<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
'<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
'<pkg:xmlData>' + `
'<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
'<w:body>' + `
'<w:p>' + `
'<w:r>' + `
'<w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr>' + `
'<w:t>The data suggests that there is more variability with successful campaigns since the variance for successful campaigns is greater than the one for unsuccessful campaigns.  For a campaign to be successful, it has to attract a large number of investors, this will likely increase the variance as a result. At the same time</w:t>' + `
'</w:r>' + `
'<w:r>' + `
'<w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr>' + `
'<w:t>,</w:t>' + `
'</w:r>' + `
'<w:r>' + `
'<w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr>' + `
'<w:t xml:space="preserve"> the data may mask other factors for example the campaigns which were not successful may have not been promoted well enough.</w:t>' + `
'</w:r>' + `
'</w:p>' + `
'</w:body>' + `
'</w:document>' + `
'</pkg:xmlData>' + `
'</pkg:part>' + `
'</pkg:package>'

$null = $targetRng.InsertXML($xmlFragment)
